# Update the "想去人数" (want-to-go headcount) column F values across the
# sheets that track event listings. Sheet3 (本地生活) has no matching rows
# in the diff, so it is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 92
$ws1.Range("F5").Value = 9944
$ws1.Range("F6").Value = 699
$ws1.Range("F9").Value = 183
$ws1.Range("F10").Value = 367
$ws1.Range("F11").Value = 447
$ws1.Range("F12").Value = 227
$ws1.Range("F13").Value = 23
$ws1.Range("F14").Value = 501
$ws1.Range("F15").Value = 12568
$ws1.Range("F16").Value = 12568
$ws1.Range("F23").Value = 395
$ws1.Range("F24").Value = 260
$ws1.Range("F27").Value = 136
$ws1.Range("F36").Value = 1073
$ws1.Range("F37").Value = 4250
$ws1.Range("F39").Value = 792
$ws1.Range("F41").Value = 58
$ws1.Range("F42").Value = 1360
$ws1.Range("F44").Value = 42
$ws1.Range("F45").Value = 475
$ws1.Range("F46").Value = 628
$ws1.Range("F47").Value = 75
$ws1.Range("F48").Value = 278

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 9

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 92
$ws4.Range("F6").Value = 9945
$ws4.Range("F7").Value = 699
$ws4.Range("F10").Value = 183
$ws4.Range("F11").Value = 367
$ws4.Range("F12").Value = 227
$ws4.Range("F13").Value = 23
$ws4.Range("F14").Value = 12568
$ws4.Range("F20").Value = 260
$ws4.Range("F23").Value = 136
$ws4.Range("F31").Value = 1073
$ws4.Range("F35").Value = 4250
$ws4.Range("F37").Value = 792
$ws4.Range("F39").Value = 58
$ws4.Range("F41").Value = 1360
$ws4.Range("F43").Value = 42
$ws4.Range("F44").Value = 475
$ws4.Range("F46").Value = 628
$ws4.Range("F48").Value = 278
